# Updated cryptos list (Price + Volume(1h) columns) to match latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.542.48"
$ws.Range("E2").Value = "  -2.24%  "

$ws.Range("D3").Value = "2.672.26"
$ws.Range("E3").Value = "  -2.90%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.99%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.36%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").Value = "  -0.41%  "

$ws.Range("D9").Value = "2.671.83"
$ws.Range("E9").Value = "  -2.87%  "

$ws.Range("E10").Value = "  +0.89%  "

$ws.Range("E11").Value = "  +1.16%  "

$ws.Range("E12").Value = "  -0.73%  "

$ws.Range("E13").Value = "  -2.55%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.85"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.78%  "

$ws.Range("D15").Value = "3.158.88"
$ws.Range("E15").Value = "  -2.88%  "

$ws.Range("E16").Value = "  -3.23%  "

$ws.Range("D17").Value = "67.471.78"
$ws.Range("E17").Value = "  -2.20%  "

$ws.Range("D18").Value = "2.646.01"
$ws.Range("E18").Value = "  -4.40%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.56%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.57%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "364.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.39%  "

$ws.Range("E22").Value = "  -4.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.80"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.71%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.29%  "

$ws.Range("E25").Value = "  +0.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "70.78"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.64%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.77%  "

$ws.Range("E29").Value = "  -4.35%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.05%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "551.87"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -8.31%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.32%  "

$ws.Range("E33").Value = "  -5.03%  "

$ws.Range("E34").Value = "  -2.60%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.131"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.07%  "

$ws.Range("E36").Value = "  +0.05%  "

$ws.Range("E37").Value = "  -5.61%  "

$ws.Range("E38").Value = "  -3.78%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "155.42"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.88%  "

$ws.Range("E40").Value = "  -3.35%  "

$ws.Range("E41").Value = "  -4.73%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.26"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.90%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.93"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.64%  "

$ws.Range("E45").Value = "  -7.70%  "

$ws.Range("E46").Value = "  -1.12%  "

$ws.Range("E47").Value = "  -6.26%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.590"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.68%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "153.67"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.59%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.86"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.88%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.70%  "

